$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.744.06"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "'2.303.73"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'271.51"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "'93.53"
$ws.Range("E6").Value = "  +5.77%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.620"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").Value = "'44.63"
$ws.Range("E10").Value = "  -3.31%  "
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Value = "'8.11"
$ws.Range("E12").Value = "  +6.91%  "
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "'2.650.76"
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("D15").Value = "'15.28"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "'0.843"
$ws.Range("E16").Value = "  +5.72%  "
$ws.Range("D17").Value = "'2.303.45"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("D18").Value = "'43.749.29"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "'0.0000106"
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("D20").Value = "'6.26"
$ws.Range("E20").Value = "  +3.57%  "
$ws.Range("D21").Value = "'71.47"
$ws.Range("E21").Value = "  +1.76%  "
$ws.Range("D22").Value = "'239.93"
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("E23").Value = "  -5.66%  "
$ws.Range("D24").Value = "'9.70"
$ws.Range("E24").Value = "  +8.75%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'11.34"
$ws.Range("E26").Value = "  +3.50%  "
$ws.Range("E27").Value = "  -3.61%  "
$ws.Range("E28").Value = "  +5.73%  "
$ws.Range("D29").Value = "'3.38"
$ws.Range("E29").Value = "  -4.59%  "
$ws.Range("D30").Value = "'38.98"
$ws.Range("E30").Value = "  -4.18%  "
$ws.Range("D31").Value = "'22.63"
$ws.Range("E31").Value = "  +9.20%  "
$ws.Range("D32").Value = "'171.53"
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("D33").Value = "'0.0899"
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("D34").Value = "'5.57"
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("D37").Value = "'4.49"
$ws.Range("E37").Value = "  +2.23%  "
$ws.Range("D38").Value = "'0.0354"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").Value = "'3.42"
$ws.Range("E39").Value = "  +2.21%  "
$ws.Range("E40").Value = "  +15.16%  "
$ws.Range("D41").Value = "'2.28"
$ws.Range("E41").Value = "  +6.59%  "
$ws.Range("D42").Value = "'12.13"
$ws.Range("E42").Value = "  -4.24%  "
$ws.Range("D43").Value = "'1.32"
$ws.Range("E43").Value = "  +15.59%  "
$ws.Range("D44").Value = "'5.44"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").Value = "'61.81"
$ws.Range("E45").Value = "  -5.89%  "
$ws.Range("E46").Value = "  +6.81%  "
$ws.Range("D47").Value = "'0.102"
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("D48").Value = "'100.12"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("D50").Value = "'2.528.70"
$ws.Range("E50").Value = "  +2.79%  "
$ws.Range("D51").Value = "'0.424"
$ws.Range("E51").Value = "  -4.11%  "
